$d = $word.ActiveDocument

# --- Step 1: remove the old _GoBack bookmark (it will be relocated) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: extend the "Los datos..." paragraph with the new sentence,
#             the relocated _GoBack bookmark, and two new paragraphs
#             (one blank, one with the Google geocoding API link). ---
$r = $d.Content
$r.Find.Execute("Los datos están desglosados por cada uno de los 28 paises de la UE.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)

$frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="62BDFB82" w14:textId="77777777" w:rsidR="00A67A8E" w:rsidRDefault="00A67A8E" w:rsidP="00F80A2A"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="300" w:lineRule="exact"/><w:ind w:left="284" w:firstLine="567"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Los datos están desglosados por cada uno de los 28 paises de la UE</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>. Cada país se ha geolocalizado con un p</w:t></w:r><w:bookmarkStart w:id="3" w:name="_GoBack"/><w:bookmarkEnd w:id="3"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>unto representativo y se ha expresado en dos campos (longitud y latitud). Para ello se ha utilizado la API de geocodificación de Google:</w:t></w:r></w:p><w:p><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="300" w:lineRule="exact"/><w:ind w:left="284" w:firstLine="567"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="300" w:lineRule="exact"/><w:ind w:left="284" w:firstLine="567"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>https://developers.google.com/maps/documentation/geocoding/intro?hl=es-419</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($frag)

Write-Host "Edit complete"
